$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.097472071647644
$ws.Range("B1").Value = 3.556003332138062
$ws.Range("C1").Value = 3.31883978843689
$ws.Range("D1").Value = 3.747977495193481
$ws.Range("E1").Value = 1.098143815994263
